# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns
# on rows 2-51 to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.579.12'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '3.609.75'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '''203.40'
$ws.Range("E5").Value = '  +4.60%  '
$ws.Range("D6").Value = '''596.60'
$ws.Range("E6").Value = '  -1.59%  '
$ws.Range("D7").Value = '''0.628'
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +6.79%  '
$ws.Range("D11").Value = '''53.99'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '''9.68'
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = '4.179.55'
$ws.Range("E14").Value = '  +2.14%  '
$ws.Range("D15").Value = '''682.72'
$ws.Range("E15").Value = '  +15.36%  '
$ws.Range("D16").Value = '70.628.51'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '''19.15'
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").Value = '''12.77'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").Value = '3.614.41'
$ws.Range("E19").Value = '  +3.22%  '
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").Value = '''0.998'
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").Value = '''18.44'
$ws.Range("E22").Value = '  +4.20%  '
$ws.Range("D23").Value = '''110.53'
$ws.Range("E23").Value = '  +7.34%  '
$ws.Range("D24").Value = '''5.27'
$ws.Range("E24").Value = '  +3.11%  '
$ws.Range("D25").Value = '''4.55'
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("D26").Value = '''3.02'
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = '''10.63'
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").Value = '''10.09'
$ws.Range("E29").Value = '  +6.23%  '
$ws.Range("D30").Value = '''34.42'
$ws.Range("E30").Value = '  +4.07%  '
$ws.Range("E31").Value = '  +4.86%  '
$ws.Range("D32").Value = '''7.18'
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("D33").Value = '''12.36'
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("D35").Value = '''63.62'
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = '0.0₃0852'
$ws.Range("E36").Value = '  +4.02%  '
$ws.Range("D37").Value = '3.861.42'
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = '''509.94'
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = '''3.01'
$ws.Range("E40").Value = '  -5.65%  '
$ws.Range("D41").Value = '''3.59'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '''36.81'
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("E44").Value = '  +3.23%  '
$ws.Range("D45").Value = '''0.0472'
$ws.Range("E45").Value = '  +5.67%  '
$ws.Range("D46").Value = '''3.08'
$ws.Range("E46").Value = '  +9.94%  '
$ws.Range("D47").Value = '''3.45'
$ws.Range("E47").Value = '  +3.45%  '
$ws.Range("E48").Value = '  +1.64%  '
$ws.Range("D49").Value = '''8.66'
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").Value = '''1.84'
$ws.Range("E51").Value = '  +24.10%  '
